# Update gh-pages to output generated at 456a3b4
#
# The upstream scrape re-ran and:
#  - dropped the "合肥·银魂主题派对only2.0" listing (2024-08-17) from the
#    "展览" and "全部类型" sheets (it shows up as unsellable / stale),
#    shifting all subsequent rows up by one and renumbering the serial
#    (column A) index accordingly;
#  - refreshed the "想去人数" (interest count, column F) for several
#    still-live listings across all three populated sheets.
#
# "本地生活" has no data rows and is left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Remove the stale "银魂主题派对only2.0" row; everything below shifts up.
$ws1.Rows(3).Delete()

# Refresh interest counts ("想去人数", column F).
$ws1.Range("F2").Value = 5729
$ws1.Range("F3").Value = 863
$ws1.Range("F4").Value = 80
$ws1.Range("F5").Value = 401

# Renumber the serial index column (A) for the shifted rows.
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4
$ws1.Range("A6").Value = 5
$ws1.Range("A7").Value = 6
$ws1.Range("A8").Value = 7

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) - no row added/removed, just count refresh.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 60
$ws2.Range("F3").Value = 24

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - mirrors "展览" plus the "演出" rows.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Remove the same stale "银魂主题派对only2.0" row.
$ws4.Rows(3).Delete()

# Refresh interest counts ("想去人数", column F).
$ws4.Range("F2").Value = 5729
$ws4.Range("F3").Value = 863
$ws4.Range("F4").Value = 80
$ws4.Range("F5").Value = 60
$ws4.Range("F6").Value = 401
$ws4.Range("F10").Value = 24

# Renumber the serial index column (A) for the shifted rows.
$ws4.Range("A3").Value = 2
$ws4.Range("A4").Value = 3
$ws4.Range("A5").Value = 4
$ws4.Range("A6").Value = 5
$ws4.Range("A7").Value = 6
$ws4.Range("A8").Value = 7
$ws4.Range("A9").Value = 8
$ws4.Range("A10").Value = 9
$ws4.Range("A11").Value = 10
